$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Single-value cell updates
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "37"
$t.Cell(5,1).Range.Text  = "0.00003"
$t.Cell(6,1).Range.Text  = "0.00074"
$t.Cell(7,1).Range.Text  = "0.00021"
$t.Cell(8,1).Range.Text  = "0.00008"
$t.Cell(9,1).Range.Text  = "0.00030"
$t.Cell(10,1).Range.Text = "0.00044"
$t.Cell(11,1).Range.Text = "0.00070"
$t.Cell(12,1).Range.Text = "0.00938"

# Cells that collapse a tab-separated run of several values down to a single value
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "65"
